$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A new daily price record was inserted as row 127 (Hortaliza / weekly->daily
# logic): every existing record from the old row 127 down to the old row 224
# shifts down by one row, and the data that falls off the bottom (old row
# 224) becomes the new row 225. Insert a whole row at 127 so Excel performs
# that shift (and grows the used range / dimension to R225) automatically.
$ws.Rows("127:127").Insert()

# Populate the newly inserted row 127 with the new record's data.
$ws.Range("A127").Value = 3
$ws.Range("B127").Value = "Femacal de La Calera"
$ws.Range("C127").Value = "Coquimbo"
$ws.Range("D127").Value = 44957
$ws.Range("E127").Value = 5
$ws.Range("F127").Value = 100112052
$ws.Range("G127").Value = "Albahaca"
$ws.Range("H127").Value = "Sin especificar"
$ws.Range("I127").Value = "Primera"
$ws.Range("J127").Value = 100
$ws.Range("K127").Value = 4500
$ws.Range("L127").Value = 5000
$ws.Range("M127").Value = 4750
$ws.Range("N127").Value = "$/docena de matas"
$ws.Range("O127").Value = "Provincia de Quillota"
$ws.Range("P127").Value = 792
$ws.Range("Q127").Value = 6
$ws.Range("R127").Value = "Hortaliza"
